# Auto-generated edit script: updates market-data snapshot values (columns H-N)
# for specific Leve rows across all 8 job sheets, per the scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4322.6665
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 4322.6665
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 4322.6665
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -5404.6665

$ws.Range("H111").Value = 15629419
$ws.Range("I111").Value = 17861342
$ws.Range("K111").Value = 53584026
$ws.Range("M111").Value = -53580959

$ws.Range("H132").Value = 1975.8
$ws.Range("I132").Value = 1634.6086
$ws.Range("K132").Value = 4903.825800000001
$ws.Range("M132").Value = -2373.825800000001

$ws.Range("H135").Value = 526941.8
$ws.Range("I135").Value = 715007.5600000001
$ws.Range("J135").Value = 357.8
$ws.Range("K135").Value = 6435068.040000001
$ws.Range("L135").Value = 3220.2
$ws.Range("M135").Value = -6432533.040000001
$ws.Range("N135").Value = -8290.200000000001

$ws.Range("H138").Value = 4150.528
$ws.Range("I138").Value = 1075.625
$ws.Range("J138").Value = 28749.75
$ws.Range("K138").Value = 3226.875
$ws.Range("L138").Value = 86249.25
$ws.Range("M138").Value = 1913.125
$ws.Range("N138").Value = -96529.25

$ws.Range("H141").Value = 2256.037
$ws.Range("I141").Value = 2150.6538
$ws.Range("K141").Value = 6451.9614
$ws.Range("M141").Value = -1271.9614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6277.7617
$ws.Range("I45").Value = 1814.8
$ws.Range("K45").Value = 1814.8
$ws.Range("M45").Value = -1437.8

$ws.Range("H61").Value = 7919.893
$ws.Range("I61").Value = 3870.923
$ws.Range("J61").Value = 11429
$ws.Range("K61").Value = 3870.923
$ws.Range("L61").Value = 11429
$ws.Range("M61").Value = -3658.923
$ws.Range("N61").Value = -11853

$ws.Range("H74").Value = 2212.9153
$ws.Range("I74").Value = 1078.9445
$ws.Range("K74").Value = 1078.9445
$ws.Range("M74").Value = -204.9445000000001

$ws.Range("H77").Value = 2212.9153
$ws.Range("I77").Value = 1078.9445
$ws.Range("K77").Value = 5394.7225
$ws.Range("M77").Value = -1026.7225

$ws.Range("H102").Value = 1581.75
$ws.Range("I102").Value = 1630.8
$ws.Range("K102").Value = 1630.8
$ws.Range("M102").Value = -8.799999999999955

$ws.Range("H110").Value = 22223678
$ws.Range("I110").Value = 1438.9
$ws.Range("K110").Value = 1438.9
$ws.Range("M110").Value = 606.0999999999999

$ws.Range("H132").Value = 3373.027
$ws.Range("I132").Value = 2004.5962
$ws.Range("J132").Value = 6607.5
$ws.Range("K132").Value = 6013.7886
$ws.Range("L132").Value = 19822.5
$ws.Range("M132").Value = -3483.7886
$ws.Range("N132").Value = -24882.5

$ws.Range("H136").Value = 7919.893
$ws.Range("I136").Value = 3870.923
$ws.Range("J136").Value = 11429
$ws.Range("K136").Value = 11612.769
$ws.Range("L136").Value = 34287
$ws.Range("M136").Value = -9062.769
$ws.Range("N136").Value = -39387

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 20834400
$ws.Range("I64").Value = 37037864
$ws.Range("J64").Value = 1377.7142
$ws.Range("K64").Value = 37037864
$ws.Range("L64").Value = 1377.7142
$ws.Range("M64").Value = -37037639
$ws.Range("N64").Value = -1827.7142

$ws.Range("H67").Value = 20834400
$ws.Range("I67").Value = 37037864
$ws.Range("J67").Value = 1377.7142
$ws.Range("K67").Value = 37037864
$ws.Range("L67").Value = 1377.7142
$ws.Range("M67").Value = -37037084
$ws.Range("N67").Value = -2937.7142

$ws.Range("H74").Value = 23552.857
$ws.Range("J74").Value = 23552.857
$ws.Range("L74").Value = 23552.857
$ws.Range("N74").Value = -25424.857

$ws.Range("H77").Value = 23552.857
$ws.Range("J77").Value = 23552.857
$ws.Range("L77").Value = 70658.571
$ws.Range("N77").Value = -80018.571

$ws.Range("H86").Value = 38503412
$ws.Range("I86").Value = 79351.08
$ws.Range("K86").Value = 79351.08
$ws.Range("M86").Value = -78228.08

$ws.Range("H89").Value = 38503412
$ws.Range("I89").Value = 79351.08
$ws.Range("K89").Value = 396755.4
$ws.Range("M89").Value = -391139.4

$ws.Range("H99").Value = 2528250.5
$ws.Range("I99").Value = 2803.5356
$ws.Range("J99").Value = 11367314
$ws.Range("K99").Value = 2803.5356
$ws.Range("L99").Value = 11367314
$ws.Range("M99").Value = -1305.5356
$ws.Range("N99").Value = -11370310

$ws.Range("H107").Value = 59216044
$ws.Range("I107").Value = 70316860
$ws.Range("J107").Value = 11663.333
$ws.Range("K107").Value = 70316860
$ws.Range("L107").Value = 11663.333
$ws.Range("M107").Value = -70314940
$ws.Range("N107").Value = -15503.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6613.773
$ws.Range("I31").Value = 2256.353
$ws.Range("J31").Value = 9357.333000000001
$ws.Range("K31").Value = 2256.353
$ws.Range("L31").Value = 9357.333000000001
$ws.Range("M31").Value = -1961.353
$ws.Range("N31").Value = -9947.333000000001

$ws.Range("H34").Value = 6613.773
$ws.Range("I34").Value = 2256.353
$ws.Range("J34").Value = 9357.333000000001
$ws.Range("K34").Value = 2256.353
$ws.Range("L34").Value = 9357.333000000001
$ws.Range("M34").Value = -2054.353
$ws.Range("N34").Value = -9761.333000000001

$ws.Range("H58").Value = 11910335
$ws.Range("J58").Value = 10845.632
$ws.Range("L58").Value = 10845.632
$ws.Range("N58").Value = -11251.632

$ws.Range("H111").Value = 89900
$ws.Range("J111").Value = 89900
$ws.Range("L111").Value = 89900
$ws.Range("N111").Value = -98080

$ws.Range("H118").Value = 95900
$ws.Range("J118").Value = 95900
$ws.Range("L118").Value = 95900
$ws.Range("N118").Value = -99214

$ws.Range("H122").Value = 1256.2916
$ws.Range("I122").Value = 839.61536
$ws.Range("J122").Value = 1748.7273
$ws.Range("K122").Value = 2518.84608
$ws.Range("L122").Value = 5246.1819
$ws.Range("M122").Value = -68.84608000000026
$ws.Range("N122").Value = -10146.1819

$ws.Range("H124").Value = 44998
$ws.Range("J124").Value = 44998
$ws.Range("L124").Value = 44998
$ws.Range("N124").Value = -49908

$ws.Range("H136").Value = 11910335
$ws.Range("J136").Value = 10845.632
$ws.Range("L136").Value = 32536.896
$ws.Range("N136").Value = -37636.896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3333912.8
$ws.Range("J12").Value = 4545827.5
$ws.Range("L12").Value = 13637482.5
$ws.Range("N12").Value = -13637828.5

$ws.Range("H56").Value = 6884
$ws.Range("I56").Value = 6884
$ws.Range("K56").Value = 6884
$ws.Range("M56").Value = -6354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 572657.9399999999
$ws.Range("J107").Value = 927.4286
$ws.Range("L107").Value = 927.4286
$ws.Range("N107").Value = -4767.4286

$ws.Range("H113").Value = 4833.1763
$ws.Range("I113").Value = 2166.375
$ws.Range("J113").Value = 7203.6665
$ws.Range("K113").Value = 2166.375
$ws.Range("L113").Value = 7203.6665
$ws.Range("M113").Value = 3.625
$ws.Range("N113").Value = -11543.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1767
$ws.Range("I16").Value = 1301
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1301
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1131
$ws.Range("N16").Value = -2340

$ws.Range("H61").Value = 3508.7058
$ws.Range("I61").Value = 1587.4348
$ws.Range("J61").Value = 7525.909
$ws.Range("K61").Value = 1587.4348
$ws.Range("L61").Value = 7525.909
$ws.Range("M61").Value = -1385.4348
$ws.Range("N61").Value = -7929.909

$ws.Range("H113").Value = 3508.7058
$ws.Range("I113").Value = 1587.4348
$ws.Range("J113").Value = 7525.909
$ws.Range("K113").Value = 1587.4348
$ws.Range("L113").Value = 7525.909
$ws.Range("M113").Value = 582.5652
$ws.Range("N113").Value = -11865.909

$ws.Range("H132").Value = 11117899
$ws.Range("I132").Value = 23812022
$ws.Range("J132").Value = 10541.583
$ws.Range("K132").Value = 71436066
$ws.Range("L132").Value = 31624.749
$ws.Range("M132").Value = -71433536
$ws.Range("N132").Value = -36684.749

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3310.36
$ws.Range("I126").Value = 2240.0908
$ws.Range("J126").Value = 4151.2856
$ws.Range("K126").Value = 6720.2724
$ws.Range("L126").Value = 12453.8568
$ws.Range("M126").Value = -4250.2724
$ws.Range("N126").Value = -17393.8568

$ws.Range("H136").Value = 25029930
$ws.Range("I136").Value = 55556430
$ws.Range("K136").Value = 166669290
$ws.Range("M136").Value = -166666740
